$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "Absent" column (H) for the days that were missing a value
# or that had not yet been marked absent (Absent = 1 when no attendance that day).
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
